# Update the "term" ValueSet metadata sheet to the 1.1.0 release:
#   - bump the recorded Version string
#   - bump the recorded Date timestamp
#   - make sure the wrap-text/vertical-top alignment that was already present
#     on the header/body cell styles is actually flagged as applied
#     (adds applyAlignment="true" to those styles, matching how Excel
#     persists alignment once it's (re)applied through the UI).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-apply the wrap-text / top alignment on the header row and the
# bordered data cells of every sheet (only the cells that already carry
# that formatting), so the alignment actually gets flagged as applied
# (applyAlignment="true") instead of being silently ignored.
$meta.Range("A1:B1").WrapText = $true
$meta.Range("A2:B14").WrapText = $true

$ffb = $wb.Worksheets.Item("Include from FFB")
$ffb.Range("A1:C1").WrapText = $true
$ffb.Range("A2:C2").WrapText = $true
$ffb.Range("A3:B4").WrapText = $true
